$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") must stay text -- force text format before writing
# so Excel does not silently reinterpret values like "99.20" or "541.54"
# as numbers (which would drop trailing zeros / reformat the string).
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D9", "D10", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D36", "D39", "D40", "D41", "D42", "D44", "D46", "D47", "D48")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '71.889.62'
$ws.Range("E2").Value = '  +4.57%  '
$ws.Range("D3").Value = '4.037.28'
$ws.Range("E3").Value = '  +4.20%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '541.54'
$ws.Range("E5").Value = '  +3.30%  '
$ws.Range("D6").Value = '153.76'
$ws.Range("E6").Value = '  +8.73%  '
$ws.Range("D7").Value = '0.694'
$ws.Range("E7").Value = '  +13.78%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = '0.763'
$ws.Range("E9").Value = '  +6.74%  '
$ws.Range("D10").Value = '0.175'
$ws.Range("E10").Value = '  +3.94%  '
$ws.Range("E11").Value = '  +2.39%  '
$ws.Range("D12").Value = '48.42'
$ws.Range("E12").Value = '  +15.82%  '
$ws.Range("D13").Value = '10.86'
$ws.Range("E13").Value = '  +3.87%  '
$ws.Range("D14").Value = '4.691.68'
$ws.Range("E14").Value = '  +4.68%  '
$ws.Range("D15").Value = '4.049.60'
$ws.Range("E15").Value = '  +3.82%  '
$ws.Range("D16").Value = '14.44'
$ws.Range("E16").Value = '  +2.06%  '
$ws.Range("D17").Value = '20.66'
$ws.Range("E17").Value = '  -3.64%  '
$ws.Range("D18").Value = '1.21'
$ws.Range("E18").Value = '  +0.99%  '
$ws.Range("E19").Value = '  -0.23%  '
$ws.Range("D20").Value = '71.821.56'
$ws.Range("E20").Value = '  +4.55%  '
$ws.Range("D21").Value = '435.28'
$ws.Range("E21").Value = '  +4.10%  '
$ws.Range("D22").Value = '99.20'
$ws.Range("E22").Value = '  +13.89%  '
$ws.Range("D23").Value = '3.58'
$ws.Range("E23").Value = '  +1.56%  '
$ws.Range("D24").Value = '4.29'
$ws.Range("E24").Value = '  +6.50%  '
$ws.Range("D25").Value = '14.69'
$ws.Range("E25").Value = '  +4.36%  '
$ws.Range("D26").Value = '11.31'
$ws.Range("E26").Value = '  -4.71%  '
$ws.Range("D27").Value = '10.94'
$ws.Range("E27").Value = '  +4.23%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = '3.70'
$ws.Range("E28").Value = '  +29.77%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '37.08'
$ws.Range("E29").Value = '  +4.11%  '
$ws.Range("B30").Value = 'LEO'
$ws.Range("C30").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D30").Value = '5.83'
$ws.Range("E30").Value = '  +2.70%  '
$ws.Range("D31").Value = '13.77'
$ws.Range("E31").Value = '  +1.99%  '
$ws.Range("E32").Value = '  +5.67%  '
$ws.Range("D33").Value = '684.66'
$ws.Range("E33").Value = '  +0.69%  '
$ws.Range("D34").Value = '6.94'
$ws.Range("E34").Value = '  +0.20%  '
$ws.Range("D35").Value = '67.37'
$ws.Range("E35").Value = '  +0.46%  '
$ws.Range("D36").Value = '43.14'
$ws.Range("E36").Value = '  +8.72%  '
$ws.Range("E37").Value = '  -2.01%  '
$ws.Range("E38").Value = '  +6.11%  '
$ws.Range("D39").Value = '0.0₃0847'
$ws.Range("E39").Value = '  -0.13%  '
$ws.Range("D40").Value = '3.45'
$ws.Range("E40").Value = '  -1.90%  '
$ws.Range("D41").Value = '3.44'
$ws.Range("E41").Value = '  +9.16%  '
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.21%  '
$ws.Range("E43").Value = '  +4.22%  '
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("E45").Value = '  +7.49%  '
$ws.Range("D46").Value = '2.74'
$ws.Range("E46").Value = '  -3.71%  '
$ws.Range("D47").Value = '3.41'
$ws.Range("E47").Value = '  +0.13%  '
$ws.Range("D48").Value = '9.52'
$ws.Range("E48").Value = '  +8.11%  '
$ws.Range("E49").Value = '  +1.56%  '
$ws.Range("E50").Value = '  +2.03%  '
$ws.Range("E51").Value = '  -2.36%  '

# Restore default (unstyled) cell style now that the text value is set,
# matching the original formatting (no explicit style on these cells).
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
